# Monthly report sample: correct the "Lot No" header labels to include a
# period ("Lot No. 1" / "Lot No. 2"), and move the sheet's active selection
# over to the second "Lot No." column (Z1), scrolling the view so column K
# is the new left-most visible column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T1").Value = "Lot No. 1"
$ws.Range("Z1").Value = "Lot No. 2"

$ws.Range("Z1").Select()
$excel.ActiveWindow.ScrollColumn = 11
